$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing quarter header labels (use .Text because .Value
# does not reliably resolve when read back in this runtime).
$textB1 = $ws.Range("B1").Text
$textC1 = $ws.Range("C1").Text
$textD1 = $ws.Range("D1").Text
$textE1 = $ws.Range("E1").Text

# Shift the quarter header labels in row 1 one column to the right
# (B1:E1 -> C1:F1), then insert the new "Sep '13" label in B1.
$ws.Range("F1").Value = $textE1
$ws.Range("E1").Value = $textD1
$ws.Range("D1").Value = $textC1
$ws.Range("C1").Value = $textB1
$ws.Range("B1").Value = "Sep '13"
